$d = $word.ActiveDocument

# --- Step 1: perform the textual replacement -------------------------------
# This touches the run holding the sentence. The runtime's Find/Replace
# normalizes (merges) every subsequent run in the paragraph that shares the
# same run properties, so after this call the remainder of the paragraph
# (through "...заявки.") collapses into a single run. We restore the
# original run boundaries (and create the new ones requested by the edit)
# explicitly afterwards.

$old = " пользователи заходят на сайт для того, чтобы максимально быстро подобрать нужное им оформление квартиры, поэтому если при "
$new = " пользователи заходят на сайт для того, чтобы максимально быстро подобрать нужную им  квартиру, поэтому если при "

$hit = $d.Content
$found = $hit.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) {
    throw "Target sentence for replacement was not found."
}

$runStart = $hit.Start

# Keep the restoration work confined to the paragraph we just edited so nothing
# elsewhere in the document is touched.
$para = $d.Range($runStart, $runStart).Paragraphs(1)
$paraEnd = $para.Range.End - 1

# --- Step 2: re-establish run boundaries ------------------------------------
# Offsets (relative to $runStart) where a new <w:r> must begin, in the text
# that now reads:
#   " пользователи ... подобрать нужн" | "ую" | " им  квартир" | "у" |
#   ", поэтому если при " | "попадании" | " на сайт ... остальн" |
#   "ые действия" | ", например ... нужно " | "1. Выбрать ... заявки."
$bounds = @(79, 81, 93, 94, 113, 122, 249, 260, 312)

foreach ($b in $bounds) {
    $spanStart = $runStart + $b
    $span = $d.Range($spanStart, $paraEnd)
    $span.Font.Bold = $true
    $span.Font.Bold = $false
}
